$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45202 = 2023-10-03)
# for every data row (rows 2-452). Update it to 45203 (2023-10-04).
$ws.Range("C2:C452").Value = 45203
